# B2 Token whitepaper update
# 1. "B2 is a deflationary..." -> "B2 Token is a deflationary..."
# 2. Contract address becomes a hyperlink to a new address (bscscan-style explorer link)
# 3. "Liquidity:" wallet bullet is renamed to "Deployer:" and loses its trailing parenthetical
# 4. A new "Tax Wallet:" bullet is added right after "Founder Reserve:"

$d = $word.ActiveDocument

# --- 1. Intro paragraph tweak -------------------------------------------------
$d.Content.Find.Execute(
    "B2 is a deflationary utility token deployed on BNB Smart Chain.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "B2 Token is a deflationary utility token deployed on BNB Smart Chain.",
    2
) | Out-Null

# --- 2. Contract address -> hyperlink -----------------------------------------
$addrFind = $d.Content.Find
$addrFind.Execute("0x95a7fa18a399d1d4498c1662f8e60918f8d46b36") | Out-Null
$addrRange = $addrFind.Parent
$addrRange.Text = "0x8fc4815EAd6d8d0ec93D6132f14c5E9eC85dfFae"
$d.Hyperlinks.Add(
    $addrRange,
    "https://bscscan.com/address/0x8fc4815EAd6d8d0ec93D6132f14c5E9eC85dfFae",
    $null,
    $null,
    "0x8fc4815EAd6d8d0ec93D6132f14c5E9eC85dfFae"
) | Out-Null

# --- 3. Liquidity wallet -> Deployer wallet -----------------------------------
# (the document's bullet values consistently lead with a non-breaking space,
# so the replacement text below re-uses that same NBSP character, not a
# plain ASCII space, to stay byte-for-byte consistent with its neighbours)
$nbsp = [char]0x00A0

$d.Content.Find.Execute(
    "Liquidity:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployer:",
    2
) | Out-Null

$d.Content.Find.Execute(
    "$($nbsp)0xc1eeba5f1d4cf00b7871a05663B890cf2C10a187 (Seeded with initial BNB/B2 pairing to provide trading liquidity)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$($nbsp)0xc1eeba5f1d4cf00b7871a05663B890cf2C10a187",
    2
) | Out-Null

# --- 4. New "Tax Wallet" bullet, inserted after "Founder Reserve" ------------
$founderIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Founder Reserve:*") {
        $founderIndex = $i
    }
}

$founderPara = $d.Paragraphs.Item($founderIndex)
$founderPara.Range.InsertParagraphAfter() | Out-Null

$newIndex = $founderIndex + 1
$newPara = $d.Paragraphs.Item($newIndex)

# Clone the Founder Reserve run/formatting (bold label + plain value runs)
# into the freshly inserted empty paragraph, then swap in the new text -
# this keeps the exact two-run (bold label / plain value) shape intact.
$newPara.Range.FormattedText = $founderPara.Range.FormattedText

$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Find.Execute(
    "Founder Reserve:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tax Wallet:",
    2
) | Out-Null

$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Find.Execute(
    "$($nbsp)0x56Cd9120B63c9a07e45Fd2ef6729BA447c85Ed74 (Set aside as a long-term incentive for project growth and alignment)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$($nbsp)0xEce40D86917b7f77D0BbEeC1F870000A98a1a27A",
    2
) | Out-Null
